$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without Excel auto-converting
# number-like strings (e.g. ""2.28"") into real numbers, and without leaving a
# lingering custom style on the cell (ClearFormats drops the style index again
# after the temporary '@' text format forced literal-text interpretation).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "44.039.63"
Set-TextValue $ws.Range("E2") "  +1.16%  "
Set-TextValue $ws.Range("D3") "2.330.95"
Set-TextValue $ws.Range("E3") "  +4.59%  "
Set-TextValue $ws.Range("E4") "  -0.14%  "
Set-TextValue $ws.Range("D5") "97.69"
Set-TextValue $ws.Range("E5") "  +6.02%  "
Set-TextValue $ws.Range("D6") "271.63"
Set-TextValue $ws.Range("E6") "  +1.02%  "
Set-TextValue $ws.Range("D7") "0.628"
Set-TextValue $ws.Range("E7") "  +0.78%  "
Set-TextValue $ws.Range("E8") "  -0.06%  "
Set-TextValue $ws.Range("D9") "0.628"
Set-TextValue $ws.Range("E9") "  +1.54%  "
Set-TextValue $ws.Range("D10") "46.39"
Set-TextValue $ws.Range("E10") "  -0.32%  "
Set-TextValue $ws.Range("D11") "0.0953"
Set-TextValue $ws.Range("E11") "  +3.37%  "
Set-TextValue $ws.Range("E12") "  -2.28%  "
Set-TextValue $ws.Range("D13") "0.106"
Set-TextValue $ws.Range("E13") "  +0.75%  "
Set-TextValue $ws.Range("D14") "2.678.04"
Set-TextValue $ws.Range("D15") "15.60"
Set-TextValue $ws.Range("E15") "  +3.48%  "
Set-TextValue $ws.Range("E16") "  +9.48%  "
Set-TextValue $ws.Range("D17") "2.331.41"
Set-TextValue $ws.Range("E17") "  +4.34%  "
Set-TextValue $ws.Range("D18") "43.920.30"
Set-TextValue $ws.Range("E18") "  +0.91%  "
Set-TextValue $ws.Range("E19") "  +5.98%  "
Set-TextValue $ws.Range("D20") "6.43"
Set-TextValue $ws.Range("E20") "  +7.19%  "
Set-TextValue $ws.Range("D21") "72.92"
Set-TextValue $ws.Range("E21") "  +3.67%  "
Set-TextValue $ws.Range("D22") "240.07"
Set-TextValue $ws.Range("E22") "  +3.12%  "
Set-TextValue $ws.Range("D23") "2.28"
Set-TextValue $ws.Range("E23") "  -1.38%  "
Set-TextValue $ws.Range("D24") "9.47"
Set-TextValue $ws.Range("E24") "  +5.90%  "
Set-TextValue $ws.Range("E25") "  -0.11%  "
Set-TextValue $ws.Range("D26") "2.53"
Set-TextValue $ws.Range("E26") "  +1.74%  "
Set-TextValue $ws.Range("D27") "11.42"
Set-TextValue $ws.Range("E27") "  +1.12%  "
Set-TextValue $ws.Range("D28") "3.47"
Set-TextValue $ws.Range("E28") "  -2.11%  "
Set-TextValue $ws.Range("D29") "2.26"
Set-TextValue $ws.Range("E29") "  -0.16%  "
Set-TextValue $ws.Range("D30") "38.42"
Set-TextValue $ws.Range("E30") "  -3.89%  "
Set-TextValue $ws.Range("D31") "22.47"
Set-TextValue $ws.Range("E31") "  +8.21%  "
Set-TextValue $ws.Range("D32") "174.50"
Set-TextValue $ws.Range("E32") "  +1.03%  "
Set-TextValue $ws.Range("D33") "0.0910"
Set-TextValue $ws.Range("E33") "  -1.50%  "
Set-TextValue $ws.Range("E34") "  +1.11%  "
Set-TextValue $ws.Range("E35") "  +3.39%  "
Set-TextValue $ws.Range("E36") "  +3.87%  "
Set-TextValue $ws.Range("E37") "  -1.33%  "
Set-TextValue $ws.Range("E38") "  +3.60%  "
Set-TextValue $ws.Range("E39") "  -4.95%  "
Set-TextValue $ws.Range("D40") "0.241"
Set-TextValue $ws.Range("E40") "  +10.88%  "
Set-TextValue $ws.Range("E41") "  +8.95%  "
Set-TextValue $ws.Range("D42") "1.37"
Set-TextValue $ws.Range("E42") "  +18.98%  "
Set-TextValue $ws.Range("D43") "12.32"
Set-TextValue $ws.Range("E43") "  -1.23%  "
Set-TextValue $ws.Range("E44") "  +10.06%  "
Set-TextValue $ws.Range("D45") "62.52"
Set-TextValue $ws.Range("E45") "  -0.63%  "
Set-TextValue $ws.Range("E46") "  +2.02%  "
Set-TextValue $ws.Range("D47") "0.103"
Set-TextValue $ws.Range("E47") "  +4.71%  "
Set-TextValue $ws.Range("D48") "100.55"
Set-TextValue $ws.Range("E48") "  +0.35%  "
Set-TextValue $ws.Range("E49") "  +1.60%  "
Set-TextValue $ws.Range("D50") "2.555.91"
Set-TextValue $ws.Range("E50") "  +4.14%  "
Set-TextValue $ws.Range("E51") "  +16.23%  "
